$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("M2").Value = 1.686881333333333
$ws.Range("N2").Value = 5.060644
$ws.Range("O2").Value = 0.3041642479870916
$ws.Range("P2").Value = 0.3041642479870916
$ws.Range("Q2").Value = 286.6275500343693
$ws.Range("R2").Value = 2579.647950309324
$ws.Range("S2").Value = 0.135084967493263
$ws.Range("T2").Value = 0.135084967493263
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("O3").Value = 0.01123918175564102
$ws.Range("P3").Value = 0.01123918175564102
$ws.Range("Q3").Value = 10.591182732124
$ws.Range("R3").Value = 95.320644589116
$ws.Range("S3").Value = 0.004991528465817829
$ws.Range("T3").Value = 0.004991528465817829
$ws.Range("G4").Value = 169.915657
$ws.Range("H4").Value = 509.746971
$ws.Range("I4").Value = 0.4441184931734509
$ws.Range("J4").Value = 0.4441184931734509
$ws.Range("M4").Value = 0.9868993333333332
$ws.Range("N4").Value = 2.960698
$ws.Range("O4").Value = 0.1779493836529276
$ws.Range("P4").Value = 0.1779493836529276
$ws.Range("Q4").Value = 167.6896486161953
$ws.Range("R4").Value = 1509.206837545758
$ws.Range("S4").Value = 0.07903061212908251
$ws.Range("T4").Value = 0.07903061212908251
$ws.Range("G5").Value = 169.915657
$ws.Range("H5").Value = 509.746971
$ws.Range("I5").Value = 0.4441184931734509
$ws.Range("J5").Value = 0.4441184931734509
$ws.Range("M5").Value = 2.809842666666667
$ws.Range("N5").Value = 8.429528000000001
$ws.Range("O5").Value = 0.5066471866043398
$ws.Range("P5").Value = 0.5066471866043397
$ws.Range("Q5").Value = 477.4362627732987
$ws.Range("R5").Value = 4296.926364959689
$ws.Range("S5").Value = 0.2250113850852876
$ws.Range("T5").Value = 0.2250113850852876
$ws.Range("I6").Value = 0.1787346690539575
$ws.Range("J6").Value = 0.1787346690539575
$ws.Range("M6").Value = 1.686881333333333
$ws.Range("N6").Value = 5.060644
$ws.Range("O6").Value = 0.3041642479870916
$ws.Range("P6").Value = 0.3041642479870916
$ws.Range("Q6").Value = 115.352729248164
$ws.Range("R6").Value = 1038.174563233476
$ws.Range("S6").Value = 0.05436469620201868
$ws.Range("T6").Value = 0.05436469620201868
$ws.Range("I7").Value = 0.1787346690539575
$ws.Range("J7").Value = 0.1787346690539575
$ws.Range("O7").Value = 0.01123918175564102
$ws.Range("P7").Value = 0.01123918175564102
$ws.Range("S7").Value = 0.002008831431531775
$ws.Range("T7").Value = 0.002008831431531774
$ws.Range("I8").Value = 0.1787346690539575
$ws.Range("J8").Value = 0.1787346690539575
$ws.Range("M8").Value = 0.9868993333333332
$ws.Range("N8").Value = 2.960698
$ws.Range("O8").Value = 0.1779493836529276
$ws.Range("P8").Value = 0.1779493836529276
$ws.Range("Q8").Value = 67.486390028538
$ws.Range("R8").Value = 607.3775102568419
$ws.Range("S8").Value = 0.03180572419556173
$ws.Range("T8").Value = 0.03180572419556173
$ws.Range("I9").Value = 0.1787346690539575
$ws.Range("J9").Value = 0.1787346690539575
$ws.Range("M9").Value = 2.809842666666667
$ws.Range("N9").Value = 8.429528000000001
$ws.Range("O9").Value = 0.5066471866043398
$ws.Range("P9").Value = 0.5066471866043397
$ws.Range("Q9").Value = 192.143344023768
$ws.Range("R9").Value = 1729.290096213912
$ws.Range("S9").Value = 0.09055541722484534
$ws.Range("T9").Value = 0.09055541722484532
$ws.Range("G10").Value = 53.27463399999999
$ws.Range("H10").Value = 159.823902
$ws.Range("I10").Value = 0.1392470275793777
$ws.Range("J10").Value = 0.1392470275793778
$ws.Range("M10").Value = 1.686881333333333
$ws.Range("N10").Value = 5.060644
$ws.Range("O10").Value = 0.3041642479870916
$ws.Range("P10").Value = 0.3041642479870916
$ws.Range("Q10").Value = 89.86798563476532
$ws.Range("R10").Value = 808.8118707128879
$ws.Range("S10").Value = 0.04235396742811923
$ws.Range("T10").Value = 0.04235396742811924
$ws.Range("G11").Value = 53.27463399999999
$ws.Range("H11").Value = 159.823902
$ws.Range("I11").Value = 0.1392470275793777
$ws.Range("J11").Value = 0.1392470275793778
$ws.Range("O11").Value = 0.01123918175564102
$ws.Range("P11").Value = 0.01123918175564102
$ws.Range("Q11").Value = 3.320714486487999
$ws.Range("R11").Value = 29.886430378392
$ws.Range("S11").Value = 0.001565022651897384
$ws.Range("T11").Value = 0.001565022651897384
$ws.Range("G12").Value = 53.27463399999999
$ws.Range("H12").Value = 159.823902
$ws.Range("I12").Value = 0.1392470275793777
$ws.Range("J12").Value = 0.1392470275793778
$ws.Range("M12").Value = 0.9868993333333332
$ws.Range("N12").Value = 2.960698
$ws.Range("O12").Value = 0.1779493836529276
$ws.Range("P12").Value = 0.1779493836529276
$ws.Range("Q12").Value = 52.57670077817732
$ws.Range("R12").Value = 473.1903070035959
$ws.Range("S12").Value = 0.02477892273325248
$ws.Range("T12").Value = 0.02477892273325248
$ws.Range("G13").Value = 53.27463399999999
$ws.Range("H13").Value = 159.823902
$ws.Range("I13").Value = 0.1392470275793777
$ws.Range("J13").Value = 0.1392470275793778
$ws.Range("M13").Value = 2.809842666666667
$ws.Range("N13").Value = 8.429528000000001
$ws.Range("O13").Value = 0.5066471866043398
$ws.Range("P13").Value = 0.5066471866043397
$ws.Range("Q13").Value = 149.6933396642507
$ws.Range("R13").Value = 1347.240056978256
$ws.Range("S13").Value = 0.07054911476610864
$ws.Range("T13").Value = 0.07054911476610864
$ws.Range("G14").Value = 91.01828266666666
$ws.Range("H14").Value = 273.054848
$ws.Range("I14").Value = 0.2378998101932138
$ws.Range("J14").Value = 0.2378998101932138
$ws.Range("M14").Value = 1.686881333333333
$ws.Range("N14").Value = 5.060644
$ws.Range("O14").Value = 0.3041642479870916
$ws.Range("P14").Value = 0.3041642479870916
$ws.Range("Q14").Value = 153.5370420224569
$ws.Range("R14").Value = 1381.833378202112
$ws.Range("S14").Value = 0.07236061686369069
$ws.Range("T14").Value = 0.0723606168636907
$ws.Range("G15").Value = 91.01828266666666
$ws.Range("H15").Value = 273.054848
$ws.Range("I15").Value = 0.2378998101932138
$ws.Range("J15").Value = 0.2378998101932138
$ws.Range("O15").Value = 0.01123918175564102
$ws.Range("P15").Value = 0.01123918175564102
$ws.Range("Q15").Value = 5.673351595178667
$ws.Range("R15").Value = 51.060164356608
$ws.Range("S15").Value = 0.00267379920639403
$ws.Range("T15").Value = 0.002673799206394029
$ws.Range("G16").Value = 91.01828266666666
$ws.Range("H16").Value = 273.054848
$ws.Range("I16").Value = 0.2378998101932138
$ws.Range("J16").Value = 0.2378998101932138
$ws.Range("M16").Value = 0.9868993333333332
$ws.Range("N16").Value = 2.960698
$ws.Range("O16").Value = 0.1779493836529276
$ws.Range("P16").Value = 0.1779493836529276
$ws.Range("Q16").Value = 89.8258824848782
$ws.Range("R16").Value = 808.432942363904
$ws.Range("S16").Value = 0.04233412459503086
$ws.Range("T16").Value = 0.04233412459503085
$ws.Range("G17").Value = 91.01828266666666
$ws.Range("H17").Value = 273.054848
$ws.Range("I17").Value = 0.2378998101932138
$ws.Range("J17").Value = 0.2378998101932138
$ws.Range("M17").Value = 2.809842666666667
$ws.Range("N17").Value = 8.429528000000001
$ws.Range("O17").Value = 0.5066471866043398
$ws.Range("P17").Value = 0.5066471866043397
$ws.Range("Q17").Value = 255.7470540835271
$ws.Range("R17").Value = 2301.723486751744
$ws.Range("S17").Value = 0.1205312695280982
$ws.Range("T17").Value = 0.1205312695280982
